# Apply the evaluation-sheet updates described by the commit
# "test code generation module - update evaluations"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("QuantitativeMetrics")

# The "Assertion validity" row note (C7) no longer carries the
# "Test pass" remark - clear it back to an empty cell.
[void]$ws.Range("C7").ClearContents()

# Updated Code BLEU score (B12) and its accompanying detail
# breakdown string (C12), reflecting a new dataflow_match_score.
$ws.Range("B12").Value = 0.3120289555429744
$ws.Range("C12").Value = "{'codebleu': 0.3120289555429744, 'ngram_match_score': 0.20225288428756, 'weighted_ngram_match_score': 0.21413536313350703, 'syntax_match_score': 0.5642857142857143, 'dataflow_match_score': 0.26744186046511625}"

# The active selection on the sheet moved from C8 to B8.
[void]$ws.Range("B8").Select()
